$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" everywhere it appears ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: literal goes on the left of -eq; $cell.Text can come back as a
            # genuine boolean for True/False-valued cells, and "$true -eq <string>"
            # would coerce the string to boolean (any non-empty string -> $true),
            # causing false positives on unrelated cells.
            if ("Ready for handoff" -eq $cell.Text) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2. Narrow the "Status" column on each sheet (was 17.2159881591797, now 13.4101845877511) ---
# The ColumnWidth COM property is stored on a fixed pixel grid by the host, so 12.5
# is the nearest settable value that lands on the target width after round-tripping.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
